# Refresh the cryptos list scrape: update Price (col D) and
# Volume(1h) (col E) for every row whose figures moved.
#
# Price cells are stored as plain text (values such as "1.799.50" use
# "." as a thousands separator, and others like "2.00"/"1.00" need
# their trailing zero preserved) so numeric-looking updates are
# written with a leading apostrophe to keep Excel from reinterpreting
# them as numbers; values that already are not valid numbers do not
# need it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.420.12'
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").Value = '1.574.49'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +1.16%  '

$ws.Range("D5").Value = "'211.05"
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("E7").Value = '  +1.09%  '

$ws.Range("D8").Value = "'46.01"
$ws.Range("E8").Value = '  +4.07%  '

$ws.Range("E9").Value = '  +2.23%  '

$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("D12").Value = "'0.0878"
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("D13").Value = '1.798.80'
$ws.Range("E13").Value = '  +0.29%  '

$ws.Range("D14").Value = '1.601.26'
$ws.Range("E14").Value = '  +2.21%  '

$ws.Range("E15").Value = '  +0.46%  '

$ws.Range("E16").Value = '  -1.35%  '

$ws.Range("D17").Value = '28.401.76'
$ws.Range("E17").Value = '  +1.97%  '

$ws.Range("D18").Value = "'62.30"
$ws.Range("E18").Value = '  -1.86%  '

$ws.Range("D19").Value = "'228.28"
$ws.Range("E19").Value = '  -1.35%  '

$ws.Range("E20").Value = '  -1.37%  '

$ws.Range("E21").Value = '  -1.58%  '

$ws.Range("E22").Value = '  +1.16%  '

$ws.Range("E23").Value = '  -4.57%  '

$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = '  -1.61%  '

$ws.Range("D25").Value = "'2.00"
$ws.Range("E25").Value = '  +3.75%  '

$ws.Range("D26").Value = "'150.59"
$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = "'14.99"
$ws.Range("E27").Value = '  -1.74%  '

$ws.Range("E28").Value = '  -1.65%  '

$ws.Range("E29").Value = '  -2.44%  '

$ws.Range("E30").Value = '  +1.20%  '

$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = '  -1.87%  '

$ws.Range("D32").Value = "'0.0464"
$ws.Range("E32").Value = '  -1.89%  '

$ws.Range("E33").Value = '  -0.87%  '

$ws.Range("E34").Value = '  -0.71%  '

$ws.Range("D35").Value = '1.390.88'
$ws.Range("E35").Value = '  -1.73%  '

$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("E38").Value = '  +2.95%  '

$ws.Range("E39").Value = '  +4.04%  '

$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("E41").Value = '  -2.03%  '

$ws.Range("E42").Value = '  +1.28%  '

$ws.Range("E43").Value = '  -1.55%  '

$ws.Range("D44").Value = "'5.61"
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("D45").Value = "'1.85"
$ws.Range("E45").Value = '  +0.32%  '

$ws.Range("D46").Value = "'0.979"
$ws.Range("E46").Value = '  +1.11%  '

$ws.Range("D47").Value = "'62.26"
$ws.Range("E47").Value = '  -2.52%  '

$ws.Range("D48").Value = '1.710.85'
$ws.Range("E48").Value = '  +0.47%  '

$ws.Range("D49").Value = "'85.64"
$ws.Range("E49").Value = '  -1.12%  '

$ws.Range("D50").Value = '0.0₆0103'
$ws.Range("E50").Value = '  +3.65%  '

$ws.Range("D51").Value = "'0.0518"
$ws.Range("E51").Value = '  -0.97%  '
